$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.758.32'
$ws.Range("E2").Value = '  -1.04%  '
$ws.Range("D3").Value = '3.422.34'
$ws.Range("E3").Value = '  -1.33%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = '573.88'
$ws.Range("E5").Value = '  -0.64%  '
$ws.Range("D6").Value = '159.14'
$ws.Range("E6").Value = '  -0.63%  '
$ws.Range("D7").Value = '0.607'
$ws.Range("E7").Value = '  +3.84%  '
$ws.Range("E8").Value = '  -0.07%  '
$ws.Range("D9").Value = '3.424.99'
$ws.Range("E9").Value = '  -1.25%  '
$ws.Range("D10").Value = '7.16'
$ws.Range("E10").Value = '  -1.49%  '
$ws.Range("E11").Value = '  -1.29%  '
$ws.Range("D12").Value = '0.441'
$ws.Range("E12").Value = '  -0.02%  '
$ws.Range("D13").Value = '4.014.87'
$ws.Range("E13").Value = '  -1.44%  '
$ws.Range("E14").Value = '  +0.10%  '
$ws.Range("E15").Value = '  -2.40%  '
$ws.Range("D16").Value = '27.68'
$ws.Range("E16").Value = '  -2.60%  '
$ws.Range("D17").Value = '64.757.66'
$ws.Range("E17").Value = '  -1.07%  '
$ws.Range("D18").Value = '3.423.23'
$ws.Range("E18").Value = '  -0.59%  '
$ws.Range("D19").Value = '6.36'
$ws.Range("E19").Value = '  -0.45%  '
$ws.Range("D20").Value = '13.90'
$ws.Range("E20").Value = '  -2.13%  '
$ws.Range("D21").Value = '380.49'
$ws.Range("E21").Value = '  -2.25%  '
$ws.Range("D22").Value = '8.01'
$ws.Range("E22").Value = '  -2.14%  '
$ws.Range("D23").Value = '0.549'
$ws.Range("E23").Value = '  +0.34%  '
$ws.Range("E24").Value = '  +0.39%  '
$ws.Range("D25").Value = '72.32'
$ws.Range("E25").Value = '  -1.76%  '
$ws.Range("E26").Value = '  -2.97%  '
$ws.Range("D27").Value = '10.15'
$ws.Range("E27").Value = '  +6.44%  '
$ws.Range("D28").Value = '0.177'
$ws.Range("E28").Value = '  -0.67%  '
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  -0.03%  '
$ws.Range("E30").Value = '  +4.28%  '
$ws.Range("D31").Value = '6.22'
$ws.Range("E31").Value = '  -3.09%  '
$ws.Range("E32").Value = '  -1.77%  '
$ws.Range("D33").Value = '23.23'
$ws.Range("E33").Value = '  -1.45%  '
$ws.Range("D34").Value = '7.10'
$ws.Range("E34").Value = '  +0.31%  '
$ws.Range("D35").Value = '1.61'
$ws.Range("E35").Value = '  +4.98%  '
$ws.Range("D36").Value = '160.46'
$ws.Range("E36").Value = '  -1.04%  '
$ws.Range("D37").Value = '1.92'
$ws.Range("E37").Value = '  -1.06%  '
$ws.Range("D38").Value = '0.0759'
$ws.Range("E38").Value = '  -1.14%  '
$ws.Range("D39").Value = '2.896.65'
$ws.Range("E39").Value = '  -4.82%  '
$ws.Range("D40").Value = '6.77'
$ws.Range("E40").Value = '  +4.33%  '
$ws.Range("D41").Value = '26.54'
$ws.Range("E41").Value = '  -1.62%  '
$ws.Range("D42").Value = '4.60'
$ws.Range("E42").Value = '  +2.42%  '
$ws.Range("D43").Value = '43.00'
$ws.Range("E43").Value = '  +0.79%  '
$ws.Range("D44").Value = '0.0316'
$ws.Range("E44").Value = '  -1.12%  '
$ws.Range("D45").Value = '0.772'
$ws.Range("E45").Value = '  -0.12%  '
$ws.Range("D46").Value = '25.92'
$ws.Range("E46").Value = '  +0.86%  '
$ws.Range("D47").Value = '319.23'
$ws.Range("E47").Value = '  +3.19%  '
$ws.Range("D48").Value = '2.25'
$ws.Range("E48").Value = '  +2.51%  '
$ws.Range("E49").Value = '  -3.87%  '
$ws.Range("D50").Value = '0.108'
$ws.Range("E50").Value = '  +1.43%  '
$ws.Range("D51").Value = '6.56'
$ws.Range("E51").Value = '  -1.58%  '
